# Scheduled market-data refresh: update cached price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across all Leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 168.6
$ws.Range("I2").Value = 148.25
$ws.Range("J2").Value = 250
$ws.Range("K2").Value = 148.25
$ws.Range("L2").Value = 250
$ws.Range("M2").Value = -35.25
$ws.Range("N2").Value = -476
$ws.Range("H6").Value = 60.583332
$ws.Range("I6").Value = 60.583332
$ws.Range("K6").Value = 181.749996
$ws.Range("M6").Value = -69.74999600000001
$ws.Range("H40").Value = 4450
$ws.Range("J40").Value = 4000
$ws.Range("L40").Value = 4000
$ws.Range("N40").Value = -4350
$ws.Range("H74").Value = 10032.429
$ws.Range("I74").Value = 9958
$ws.Range("K74").Value = 9958
$ws.Range("M74").Value = -9022
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H77").Value = 10032.429
$ws.Range("I77").Value = 9958
$ws.Range("K77").Value = 49790
$ws.Range("M77").Value = -45110
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H107").Value = 566.53845
$ws.Range("I107").Value = 570.7273
$ws.Range("K107").Value = 570.7273
$ws.Range("M107").Value = 1349.2727
$ws.Range("H115").Value = 495.83334
$ws.Range("I115").Value = 495.83334
$ws.Range("K115").Value = 1487.50002
$ws.Range("M115").Value = 79.49998000000005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6757.2744
$ws.Range("I32").Value = 5899.3267
$ws.Range("K32").Value = 5899.3267
$ws.Range("M32").Value = -5612.3267
$ws.Range("H95").Value = 34306.832
$ws.Range("J95").Value = 33913
$ws.Range("L95").Value = 33913
$ws.Range("N95").Value = -39405
$ws.Range("H97").Value = 865.15
$ws.Range("I97").Value = 779.1579
$ws.Range("K97").Value = 779.1579
$ws.Range("M97").Value = -283.1579
$ws.Range("H110").Value = 808.1
$ws.Range("I110").Value = 808.1
$ws.Range("K110").Value = 808.1
$ws.Range("M110").Value = 1236.9
$ws.Range("H132").Value = 4069.4473
$ws.Range("I132").Value = 4175.4243
$ws.Range("J132").Value = 3370
$ws.Range("K132").Value = 12526.2729
$ws.Range("L132").Value = 10110
$ws.Range("M132").Value = -9996.2729
$ws.Range("N132").Value = -15170

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 11144078
$ws.Range("I86").Value = 36304.777
$ws.Range("K86").Value = 36304.777
$ws.Range("M86").Value = -35181.777
$ws.Range("H89").Value = 11144078
$ws.Range("I89").Value = 36304.777
$ws.Range("K89").Value = 181523.885
$ws.Range("M89").Value = -175907.885

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 45012.75
$ws.Range("I50").Value = 25526
$ws.Range("J50").Value = 64499.5
$ws.Range("K50").Value = 25526
$ws.Range("L50").Value = 64499.5
$ws.Range("M50").Value = -24901
$ws.Range("N50").Value = -65749.5
$ws.Range("H92").Value = 600.5
$ws.Range("J92").Value = 600.5
$ws.Range("L92").Value = 600.5
$ws.Range("N92").Value = -5592.5
$ws.Range("H107").Value = 600
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 600
$ws.Range("N107").Value = -4440
$ws.Range("M107").ClearContents()
$ws.Range("H112").Value = 54999
$ws.Range("J112").Value = 54999
$ws.Range("L112").Value = 54999
$ws.Range("N112").Value = -57953
$ws.Range("H132").Value = 6671827
$ws.Range("I132").Value = 6901856
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 20705568
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -20703038
$ws.Range("N132").Value = -8060
$ws.Range("H134").Value = 1506.8966
$ws.Range("I134").Value = 1291.0952
$ws.Range("K134").Value = 3873.2856
$ws.Range("M134").Value = -1338.2856
$ws.Range("H141").Value = 318287.53
$ws.Range("J141").Value = 318287.53
$ws.Range("L141").Value = 318287.53
$ws.Range("N141").Value = -328647.53

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1794.125
$ws.Range("I11").Value = 1794.125
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 5382.375
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -5242.375
$ws.Range("N11").ClearContents()
$ws.Range("H16").Value = 1766.6666
$ws.Range("I16").Value = 1800
$ws.Range("K16").Value = 5400
$ws.Range("M16").Value = -5227
$ws.Range("H117").Value = 939.3333
$ws.Range("I117").Value = 852.75
$ws.Range("J117").Value = 1632
$ws.Range("K117").Value = 2558.25
$ws.Range("L117").Value = 4896
$ws.Range("M117").Value = 883.75
$ws.Range("N117").Value = -11780
$ws.Range("H129").Value = 2190.6
$ws.Range("J129").Value = 2469.8333
$ws.Range("L129").Value = 7409.499899999999
$ws.Range("N129").Value = -17409.4999
$ws.Range("H131").Value = 7000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 7000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 21000
$ws.Range("N131").Value = -31080
$ws.Range("M131").ClearContents()
$ws.Range("H132").Value = 939.6667
$ws.Range("I132").Value = 939.6667
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8457.0003
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5927.0003
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 5000731.5
$ws.Range("I14").Value = 6000578
$ws.Range("J14").Value = 1500
$ws.Range("K14").Value = 6000578
$ws.Range("L14").Value = 1500
$ws.Range("M14").Value = -6000410
$ws.Range("N14").Value = -1836
$ws.Range("H113").Value = 3639.8
$ws.Range("I113").Value = 3325
$ws.Range("K113").Value = 3325
$ws.Range("M113").Value = -1155

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 29418228
$ws.Range("I40").Value = 43483930
$ws.Range("J40").Value = 8126.8184
$ws.Range("K40").Value = 43483930
$ws.Range("L40").Value = 8126.8184
$ws.Range("M40").Value = -43483794
$ws.Range("N40").Value = -8398.8184
$ws.Range("H46").Value = 4749.1665
$ws.Range("I46").Value = 4697.5
$ws.Range("K46").Value = 4697.5
$ws.Range("M46").Value = -4509.5
$ws.Range("H93").Value = 749.75
$ws.Range("I93").Value = 749.6667
$ws.Range("K93").Value = 749.6667
$ws.Range("M93").Value = 498.3333
$ws.Range("H132").Value = 6466.8
$ws.Range("I132").Value = 4514.5713
$ws.Range("K132").Value = 13543.7139
$ws.Range("M132").Value = -11013.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 50246
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 50246
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H103").Value = 35314.285
$ws.Range("J103").Value = 35314.285
$ws.Range("L103").Value = 35314.285
$ws.Range("N103").Value = -37658.285
$ws.Range("H113").Value = 800.1429000000001
$ws.Range("I113").Value = 803.8
$ws.Range("J113").Value = 791
$ws.Range("K113").Value = 2411.4
$ws.Range("L113").Value = 2373
$ws.Range("M113").Value = -241.3999999999996
$ws.Range("N113").Value = -6713
$ws.Range("H126").Value = 2779676.8
$ws.Range("I126").Value = 4631296.5
$ws.Range("J126").Value = 2247.6667
$ws.Range("K126").Value = 13893889.5
$ws.Range("L126").Value = 6743.000100000001
$ws.Range("M126").Value = -13891419.5
$ws.Range("N126").Value = -11683.0001
$ws.Range("H132").Value = 333334270
$ws.Range("I132").Value = 800
$ws.Range("J132").Value = 500001000
$ws.Range("K132").Value = 2400
$ws.Range("L132").Value = 1500003000
$ws.Range("M132").Value = 130
$ws.Range("N132").Value = -1500008060
